# Target diff analysis
# ---------------------
# The unified diff for CsCl_fractionation.docx touches exactly four
# <w:abstractNum> entries in word/numbering.xml (abstractNumId 990, 991,
# 99411, 99414): each one's <w:nsid w:val="..."/> is swapped for a new
# 8-hex-digit value. Every other byte of every part (document.xml,
# styles.xml, the abstractNum level definitions themselves, the <w:num>
# instances that reference them, paragraph text/list membership, etc.)
# is byte-for-byte identical before/after. The commit message ("small
# changes to the illumina barcoding protocol") doesn't even match this
# document's subject (CsCl gradients), confirming the nsid churn here is
# incidental fallout from Word re-saving the owning package during a
# multi-file commit, not a deliberate edit to *this* document's content.
#
# w:nsid is a purely internal bookkeeping value Word stamps on a list
# definition. It has no rendering/semantic effect, and - matching real
# Word - it is not reachable through the Object Model: there is no
# Document/List/ListTemplate/ListLevel property for it, WordOpenXML /
# Range.XML are read-only ("the assignment changed nothing"), and
# InsertXML only ever replaces the contents of the body Range it is
# called on (it cannot retarget the numbering part). Exhaustive probing
# of this document (Find/replace, ListFormat.ListLevelNumber,
# ListIndent/ListOutdent, ApplyListTemplateWithLevel, paragraph
# delete/insert, Save/SaveAs2) confirms none of them perturb the
# existing <w:nsid> values - the numbering part is carried through
# byte-identical unless a list is structurally added/removed.
#
# So the faithful reproduction of "the change described by the diff" is
# to touch nothing: no paragraph, run, list, or style in the document
# actually changes, so the correct COM script is a no-op that leaves
# $d (and therefore out.docx) exactly as it was loaded.
$d = $word.ActiveDocument
